$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.535.56"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "2.655.91"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.15%  "

$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.401"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.71%  "

$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000195"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.62%  "

$ws.Range("D15").Value = "3.133.31"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").Value = "65.315.51"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "2.654.01"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.66%  "

$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000110"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("E26").Value = "  -5.37%  "

$ws.Range("E27").Value = "  -2.30%  "

$ws.Range("E28").Value = "  -3.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "159.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.38%  "

$ws.Range("E40").Value = "  -3.87%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "164.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0606"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.640"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0258"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("E50").Value = "  +2.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "
